$wb = $excel.ActiveWorkbook

# ---- Sheet "展览": update "想去人数" (F column) counts ----
$wsExpo = $wb.Worksheets.Item('展览')
$wsExpo.Range('F4').Value = 8056
$wsExpo.Range('F5').Value = 100
$wsExpo.Range('F9').Value = 116
$wsExpo.Range('F10').Value = 484
$wsExpo.Range('F11').Value = 168
$wsExpo.Range('F13').Value = 460
$wsExpo.Range('F16').Value = 31
$wsExpo.Range('F17').Value = 5966
$wsExpo.Range('F19').Value = 283
$wsExpo.Range('F20').Value = 2027
$wsExpo.Range('F21').Value = 50
$wsExpo.Range('F22').Value = 77
$wsExpo.Range('F24').Value = 420

# ---- Sheet "演出": insert new row 2 for the Yolo Fes event ----
$wsShow = $wb.Worksheets.Item('演出')
$wsShow.Rows('2:2').Insert()
$wsShow.Range('A2').Value = 1
$wsShow.Range('B2').NumberFormat = "@"
$wsShow.Range('B2').Value = '2024-07-26'
$wsShow.Range('C2').Value = '合肥·Yolo Fes永乐庆典Vol.3·少女偶像联合演出DAY1&DAY3'
$wsShow.Range('D2').Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$wsShow.Range('E2').Value = '2024.07.26 18:00-07.28 23:59'
$wsShow.Range('F2').Value = 3
$wsShow.Range('G2').Value = 128
$wsShow.Range('H2').Value = 'https://show.bilibili.com/platform/detail.html?id=89514'
$wsShow.Range('I2').Value = '//i2.hdslb.com/bfs/openplatform/202407/aMtLMGR31721289854139.jpeg'
# Insert() shifts the two existing rows down but keeps their old serial
# numbers in column A; renumber them sequentially (2, 3) like the other sheets.
$wsShow.Range('A3').Value = 2
$wsShow.Range('A4').Value = 3

# ---- Sheet "全部类型": same F-count updates for rows 4 & 5, plus insert the new row at position 9 ----
$wsAll = $wb.Worksheets.Item('全部类型')
$wsAll.Range('F4').Value = 8056
$wsAll.Range('F5').Value = 100
$wsAll.Rows('9:9').Insert()
$wsAll.Range('A9').Value = 8
$wsAll.Range('B9').NumberFormat = "@"
$wsAll.Range('B9').Value = '2024-07-26'
$wsAll.Range('C9').Value = '合肥·Yolo Fes永乐庆典Vol.3·少女偶像联合演出DAY1&DAY3'
$wsAll.Range('D9').Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$wsAll.Range('E9').Value = '2024.07.26 18:00-07.28 23:59'
$wsAll.Range('F9').Value = 3
$wsAll.Range('G9').Value = 128
$wsAll.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=89514'
$wsAll.Range('I9').Value = '//i2.hdslb.com/bfs/openplatform/202407/aMtLMGR31721289854139.jpeg'
# Renumber the rows pushed down by the insert (their old column-A serial
# numbers are now off by one).
for ($r = 10; $r -le 27; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}
# The "环形宇宙Plus" row (id=88650) ends up one row further down after the
# insertion above (old row22 -> new row23); per the target diff its F value
# in this sheet becomes 2028 (distinct from the 2027 used on the "展览" sheet).
$wsAll.Range('F23').Value = 2028
